# 🔄 MAJ automatique BRVM via GitHub Actions
# Refresh the "Recommandations" and "Top_YTD" sheets with the latest
# BRVM market-data snapshot (sector/equity ordering, day counts, variation
# percentages, recommendation + strategy labels).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": refresh sector + equity rows (rows 2-35) ---
$ws1.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 5
$ws1.Range("D2").Value = 902.04
$ws1.Range("E2").Value = 179.02
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

$ws1.Range("A3").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 5
$ws1.Range("D3").Value = 770.13
$ws1.Range("E3").Value = 159.73
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

$ws1.Range("A4").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 5
$ws1.Range("D4").Value = 769.54
$ws1.Range("E4").Value = 155.74
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

$ws1.Range("A5").Value = "BRVM-PRESTIGE"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 5
$ws1.Range("D5").Value = 737.33
$ws1.Range("E5").Value = 150.13
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

$ws1.Range("A6").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 5
$ws1.Range("D6").Value = 588.61
$ws1.Range("E6").Value = 118.42
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

$ws1.Range("A7").Value = "BRVM - ENERGIE"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 5
$ws1.Range("D7").Value = 584.43
$ws1.Range("E7").Value = 119.31
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

$ws1.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 5
$ws1.Range("D8").Value = 475.74
$ws1.Range("E8").Value = 97.48
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

$ws1.Range("A9").Value = "BRVM-PRINCIPAL"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 236.03
$ws1.Range("E9").Value = 236.03
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"

$ws1.Range("A10").Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 233.7
$ws1.Range("E10").Value = 233.7
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

$ws1.Range("A11").Value = "BRVM – COMPOSITE TOTAL RETURN"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 139.43
$ws1.Range("E11").Value = 139.43
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"

$ws1.Range("A12").Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Range("B12").Value = 4
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 21.18
$ws1.Range("E12").Value = 4.88
$ws1.Range("F12").Value = "🟢 Achat"
$ws1.Range("G12").Value = "✅ Renforcer"

$ws1.Range("A13").Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Range("B13").Value = 2
$ws1.Range("C13").Value = 0
$ws1.Range("D13").Value = 14.68
$ws1.Range("E13").Value = 7.43
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"

$ws1.Range("A14").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B14").Value = 2
$ws1.Range("C14").Value = 1
$ws1.Range("D14").Value = 11.76
$ws1.Range("E14").Value = -1.98
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "👀 À surveiller"

$ws1.Range("A15").Value = "SMB CI (SMBC)"
$ws1.Range("B15").Value = 2
$ws1.Range("C15").Value = 0
$ws1.Range("D15").Value = 10.86
$ws1.Range("E15").Value = 7.41
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"

$ws1.Range("A16").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B16").Value = 3
$ws1.Range("C16").Value = 1
$ws1.Range("D16").Value = 6.86
$ws1.Range("E16").Value = -4.34
$ws1.Range("F16").Value = "🟢 Achat"
$ws1.Range("G16").Value = "✅ Renforcer"

$ws1.Range("A17").Value = "SICOR CI (SICC)"
$ws1.Range("B17").Value = 1
$ws1.Range("C17").Value = 0
$ws1.Range("D17").Value = 5.97
$ws1.Range("E17").Value = 5.97
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"

$ws1.Range("A18").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Range("B18").Value = 1
$ws1.Range("C18").Value = 0
$ws1.Range("D18").Value = 5.17
$ws1.Range("E18").Value = 5.17
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "➖ Neutre"

$ws1.Range("A19").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("B19").Value = 1
$ws1.Range("C19").Value = 0
$ws1.Range("D19").Value = 3.74
$ws1.Range("E19").Value = 3.74
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "➖ Neutre"

$ws1.Range("A20").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B20").Value = 1
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = 2.81
$ws1.Range("E20").Value = 4.39
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "👀 À surveiller"

$ws1.Range("A21").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("B21").Value = 1
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = 2.07
$ws1.Range("E21").Value = -2.93
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "👀 À surveiller"

$ws1.Range("A22").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B22").Value = 1
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = 2.01
$ws1.Range("E22").Value = -1.67
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "👀 À surveiller"

$ws1.Range("A23").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B23").Value = 2
$ws1.Range("C23").Value = 2
$ws1.Range("D23").Value = 0.02
$ws1.Range("E23").Value = 5.61
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "👀 À surveiller"

$ws1.Range("A24").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("B24").Value = 1
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = 0.01
$ws1.Range("E24").Value = 3.13
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "👀 À surveiller"

$ws1.Range("A25").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B25").Value = 1
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = -1.12
$ws1.Range("E25").Value = 6.36
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "👀 À surveiller"

$ws1.Range("A26").Value = "SAPH CI (SPHC)"
$ws1.Range("B26").Value = 0
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = -1.31
$ws1.Range("E26").Value = -1.31
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "➖ Neutre"

$ws1.Range("A27").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 2
$ws1.Range("D27").Value = -1.32
$ws1.Range("E27").Value = -1.38
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "👀 À surveiller"

$ws1.Range("A28").Value = "CIE CI (CIEC)"
$ws1.Range("B28").Value = 0
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = -1.85
$ws1.Range("E28").Value = -1.85
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "➖ Neutre"

$ws1.Range("A29").Value = "SICABLE CI (CABC)"
$ws1.Range("B29").Value = 0
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = -1.99
$ws1.Range("E29").Value = -1.99
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"

$ws1.Range("A30").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("B30").Value = 0
$ws1.Range("C30").Value = 1
$ws1.Range("D30").Value = -2.08
$ws1.Range("E30").Value = -2.08
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "➖ Neutre"

$ws1.Range("A31").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 2
$ws1.Range("D31").Value = -3.99
$ws1.Range("E31").Value = -4.17
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "👀 À surveiller"

$ws1.Range("A32").Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Range("B32").Value = 0
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = -5.82
$ws1.Range("E32").Value = -5.82
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "➖ Neutre"

$ws1.Range("A33").Value = "NESTLE CI (NTLC)"
$ws1.Range("B33").Value = 0
$ws1.Range("C33").Value = 2
$ws1.Range("D33").Value = -6.06
$ws1.Range("E33").Value = -3.47
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "➖ Neutre"

$ws1.Range("A34").Value = "SUCRIVOIRE (SCRC)"
$ws1.Range("B34").Value = 0
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = -6.32
$ws1.Range("E34").Value = -6.32
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "➖ Neutre"

$ws1.Range("A35").Value = "BERNABE CI (BNBC)"
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 4
$ws1.Range("D35").Value = -16.61
$ws1.Range("E35").Value = -3.45
$ws1.Range("F35").Value = "🔴 Vente"
$ws1.Range("G35").Value = "⚠️ Risque de décrochage"

# --- Sheet "Top_YTD": refresh YTD progression values (rows 2-8 + name swap rows 3-4) ---
$ws2.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Range("B2").Value = 17235.7

$ws2.Range("A3").Value = "BRVM - INDUSTRIELS"
$ws2.Range("B3").Value = 10467.6

$ws2.Range("A4").Value = "BRVM - SERVICES FINANCIERS"
$ws2.Range("B4").Value = 10451.99

$ws2.Range("A5").Value = "BRVM-PRESTIGE"
$ws2.Range("B5").Value = 9179.37

$ws2.Range("A6").Value = "BRVM - SERVICES PUBLICS"
$ws2.Range("B6").Value = 4792.18

$ws2.Range("A7").Value = "BRVM - ENERGIE"
$ws2.Range("B7").Value = 4698.6

$ws2.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Range("B8").Value = 2729.92

$ws2.Range("A9").Value = "BRVM-PRINCIPAL"
$ws2.Range("B9").Value = 236.03

$ws2.Range("A10").Value = "BRVM - CONSOMMATION DE BASE"
$ws2.Range("B10").Value = 233.7

$ws2.Range("A11").Value = "BRVM – COMPOSITE TOTAL RETURN"
$ws2.Range("B11").Value = 139.43

